$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = '63.210.14'
$ws.Range("E2").Value = '  +2.60%  '
$ws.Range("D3").Value = '3.468.43'
$ws.Range("E3").Value = '  +2.21%  '
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue "D5" '579.11'
$ws.Range("E5").Value = '  +0.36%  '
Set-TextValue "D6" '147.96'
$ws.Range("E6").Value = '  +3.42%  '
$ws.Range("D7").Value = '3.467.82'
$ws.Range("E7").Value = '  +2.22%  '
Set-TextValue "D9" '0.480'
$ws.Range("E9").Value = '  +1.45%  '
Set-TextValue "D10" '7.68'
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("E11").Value = '  +2.06%  '
Set-TextValue "D12" '0.406'
$ws.Range("E12").Value = '  +5.34%  '
$ws.Range("D13").Value = '4.064.31'
$ws.Range("E13").Value = '  +2.31%  '
Set-TextValue "D14" '29.71'
$ws.Range("E14").Value = '  +6.18%  '
$ws.Range("E15").Value = '  +2.73%  '
$ws.Range("D16").Value = '3.460.06'
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("E17").Value = '  +1.04%  '
$ws.Range("D18").Value = '63.192.93'
$ws.Range("E18").Value = '  +2.52%  '
$ws.Range("E19").Value = '  +3.47%  '
Set-TextValue "D20" '14.42'
$ws.Range("E20").Value = '  +5.03%  '
Set-TextValue "D21" '9.29'
$ws.Range("E21").Value = '  +1.78%  '
Set-TextValue "D22" '388.81'
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("E23").Value = '  +1.99%  '
Set-TextValue "D24" '74.85'
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").Value = '3.611.07'
$ws.Range("E26").Value = '  +2.27%  '
$ws.Range("E27").Value = '  +1.97%  '
$ws.Range("E28").Value = '  -1.36%  '
$ws.Range("E29").Value = '  +3.20%  '
$ws.Range("E30").Value = '  +0.00%  '
Set-TextValue "D31" '8.18'
$ws.Range("E31").Value = '  +2.40%  '
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("E33").Value = '  +0.06%  '
Set-TextValue "D34" '1.39'
$ws.Range("E34").Value = '  -0.62%  '
Set-TextValue "D35" '23.65'
$ws.Range("E35").Value = '  +1.22%  '
Set-TextValue "D36" '5.33'
$ws.Range("E36").Value = '  +3.91%  '
Set-TextValue "D37" '7.10'
$ws.Range("E37").Value = '  +2.48%  '
Set-TextValue "D38" '31.99'
$ws.Range("E38").Value = '  +16.41%  '
$ws.Range("D41").Value = '3.506.93'
$ws.Range("E41").Value = '  +2.41%  '
Set-TextValue "D42" '0.0763'
$ws.Range("E42").Value = '  +0.68%  '
Set-TextValue "D43" '0.799'
$ws.Range("E43").Value = '  +1.99%  '
$ws.Range("E44").Value = '  +4.76%  '
Set-TextValue "D47" '4.42'
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("D48").Value = '2.626.77'
$ws.Range("E48").Value = '  +6.09%  '
Set-TextValue "D49" '2.29'
$ws.Range("E49").Value = '  +12.51%  '
Set-TextValue "D50" '23.10'
$ws.Range("E50").Value = '  +1.49%  '
$ws.Range("E51").Value = '  +2.12%  '

# Row 39/40 swap: ImmutableX <-> Monero
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D39" '169.97'
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D40" '1.56'
$ws.Range("E40").Value = '  +6.24%  '

# Row 45/46 swap: OKB <-> Stacks
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D45" '1.73'
$ws.Range("E45").Value = '  +3.85%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D46" '42.30'
$ws.Range("E46").Value = '  -0.22%  '
